$wb = $excel.ActiveWorkbook

# "compounds" sheet: update NCI Thesaurus source_version (row 3, column E) to "24.10d"
$wsCompounds = $wb.Worksheets.Item("compounds")
$wsCompounds.Range("E3").Value = "24.10d"

# Move selection/active cell to E4 on the compounds sheet and make it the active tab
$wsCompounds.Activate()
$wsCompounds.Range("E4").Select()
